$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of test data (row 11), mirroring the style (wrapText) used by rows 2-9
$ws.Range("A11").Value = "Animate The Deal"
$ws.Range("B11").Value = "On Click of the deal button a card back image should move to the player (or npc) start location."
$ws.Range("C11").Value = "Problem that it moved in the wrong direction"
$ws.Range("D11").Value = "Changed the place for the two cards to move to"
$ws.Range("E11").Value = "Now Works"

$ws.Range("A11:E11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 75

# Update the view: scroll to show the new row, and select E11
$ws.Range("E11").Select()
$excel.ActiveWindow.ScrollRow = 7
